$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8781003333333334
$ws.Range("H2").Value = 2.634301
$ws.Range("I2").Value = 0.1010434633250494
$ws.Range("J2").Value = 0.1010434633250494
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2262196666666666
$ws.Range("N2").Value = 0.6786589999999999
$ws.Range("O2").Value = 0.03145179203784564
$ws.Range("P2").Value = 0.03145179203784564
$ws.Range("Q2").Value = 0.1986435647065555
$ws.Range("R2").Value = 1.787792082359
$ws.Range("S2").Value = 0.003177997995283137
$ws.Range("T2").Value = 0.003177997995283137
$ws.Range("G3").Value = 0.8781003333333334
$ws.Range("H3").Value = 2.634301
$ws.Range("I3").Value = 0.1010434633250494
$ws.Range("J3").Value = 0.1010434633250494
$ws.Range("O3").Value = 0.9636438974901603
$ws.Range("P3").Value = 0.9636438974901604
$ws.Range("Q3").Value = 6.086192439363334
$ws.Range("R3").Value = 54.77573195427001
$ws.Range("S3").Value = 0.0973699168144547
$ws.Range("T3").Value = 0.09736991681445471
$ws.Range("G4").Value = 0.8781003333333334
$ws.Range("H4").Value = 2.634301
$ws.Range("I4").Value = 0.1010434633250494
$ws.Range("J4").Value = 0.1010434633250494
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03527466666666667
$ws.Range("N4").Value = 0.105824
$ws.Range("O4").Value = 0.004904310471994002
$ws.Range("P4").Value = 0.004904310471994003
$ws.Range("Q4").Value = 0.03097469655822222
$ws.Range("R4").Value = 0.278772269024
$ws.Range("S4").Value = 0.0004955485153115818
$ws.Range("T4").Value = 0.0004955485153115819
$ws.Range("I5").Value = 0.01326751606355713
$ws.Range("J5").Value = 0.01326751606355713
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2262196666666666
$ws.Range("N5").Value = 0.6786589999999999
$ws.Range("O5").Value = 0.03145179203784564
$ws.Range("P5").Value = 0.03145179203784564
$ws.Range("Q5").Value = 0.02608290134699999
$ws.Range("R5").Value = 0.234746112123
$ws.Range("S5").Value = 0.0004172871560897753
$ws.Range("T5").Value = 0.0004172871560897753
$ws.Range("I6").Value = 0.01326751606355713
$ws.Range("J6").Value = 0.01326751606355713
$ws.Range("O6").Value = 0.9636438974901603
$ws.Range("P6").Value = 0.9636438974901604
$ws.Range("S6").Value = 0.0127851608894995
$ws.Range("T6").Value = 0.01278516088949951
$ws.Range("I7").Value = 0.01326751606355713
$ws.Range("J7").Value = 0.01326751606355713
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03527466666666667
$ws.Range("N7").Value = 0.105824
$ws.Range("O7").Value = 0.004904310471994002
$ws.Range("P7").Value = 0.004904310471994003
$ws.Range("Q7").Value = 0.004067133792
$ws.Range("R7").Value = 0.036604204128
$ws.Range("S7").Value = 0.00006506801796785188
$ws.Range("T7").Value = 0.00006506801796785189
$ws.Range("G8").Value = 3.520787
$ws.Range("H8").Value = 10.562361
$ws.Range("I8").Value = 0.4051387963370292
$ws.Range("J8").Value = 0.4051387963370292
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2262196666666666
$ws.Range("N8").Value = 0.6786589999999999
$ws.Range("O8").Value = 0.03145179203784564
$ws.Range("P8").Value = 0.03145179203784564
$ws.Range("Q8").Value = 0.7964712615443331
$ws.Range("R8").Value = 7.168241353898998
$ws.Range("S8").Value = 0.01274234116885534
$ws.Range("T8").Value = 0.01274234116885534
$ws.Range("G9").Value = 3.520787
$ws.Range("H9").Value = 10.562361
$ws.Range("I9").Value = 0.4051387963370292
$ws.Range("J9").Value = 0.4051387963370292
$ws.Range("O9").Value = 0.9636438974901603
$ws.Range("P9").Value = 0.9636438974901604
$ws.Range("Q9").Value = 24.40289156783
$ws.Range("R9").Value = 219.62602411047
$ws.Range("S9").Value = 0.390409528726687
$ws.Range("T9").Value = 0.3904095287266871
$ws.Range("G10").Value = 3.520787
$ws.Range("H10").Value = 10.562361
$ws.Range("I10").Value = 0.4051387963370292
$ws.Range("J10").Value = 0.4051387963370292
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.03527466666666667
$ws.Range("N10").Value = 0.105824
$ws.Range("O10").Value = 0.004904310471994002
$ws.Range("P10").Value = 0.004904310471994003
$ws.Range("Q10").Value = 0.1241945878293333
$ws.Range("R10").Value = 1.117751290464
$ws.Range("S10").Value = 0.001986926441486737
$ws.Range("T10").Value = 0.001986926441486738
$ws.Range("G11").Value = 0.0464
$ws.Range("H11").Value = 0.1392
$ws.Range("I11").Value = 0.005339272199663925
$ws.Range("J11").Value = 0.005339272199663925
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.2262196666666666
$ws.Range("N11").Value = 0.6786589999999999
$ws.Range("O11").Value = 0.03145179203784564
$ws.Range("P11").Value = 0.03145179203784564
$ws.Range("Q11").Value = 0.01049659253333333
$ws.Range("R11").Value = 0.09446933279999999
$ws.Range("S11").Value = 0.0001679296788572804
$ws.Range("T11").Value = 0.0001679296788572804
$ws.Range("G12").Value = 0.0464
$ws.Range("H12").Value = 0.1392
$ws.Range("I12").Value = 0.005339272199663925
$ws.Range("J12").Value = 0.005339272199663925
$ws.Range("O12").Value = 0.9636438974901603
$ws.Range("P12").Value = 0.9636438974901604
$ws.Range("Q12").Value = 0.321602576
$ws.Range("R12").Value = 2.894423184
$ws.Range("S12").Value = 0.005145157072245006
$ws.Range("T12").Value = 0.005145157072245007
$ws.Range("G13").Value = 0.0464
$ws.Range("H13").Value = 0.1392
$ws.Range("I13").Value = 0.005339272199663925
$ws.Range("J13").Value = 0.005339272199663925
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.03527466666666667
$ws.Range("N13").Value = 0.105824
$ws.Range("O13").Value = 0.004904310471994002
$ws.Range("P13").Value = 0.004904310471994003
$ws.Range("Q13").Value = 0.001636744533333333
$ws.Range("R13").Value = 0.0147307008
$ws.Range("S13").Value = 0.00002618544856163824
$ws.Range("T13").Value = 0.00002618544856163824
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.41099
$ws.Range("H14").Value = 1.23297
$ws.Range("I14").Value = 0.0472928336495663
$ws.Range("J14").Value = 0.0472928336495663
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.2262196666666666
$ws.Range("N14").Value = 0.6786589999999999
$ws.Range("O14").Value = 0.03145179203784564
$ws.Range("P14").Value = 0.03145179203784564
$ws.Range("Q14").Value = 0.09297402080333331
$ws.Range("R14").Value = 0.8367661872299998
$ws.Range("S14").Value = 0.001487444368826588
$ws.Range("T14").Value = 0.001487444368826588
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.41099
$ws.Range("H15").Value = 1.23297
$ws.Range("I15").Value = 0.0472928336495663
$ws.Range("J15").Value = 0.0472928336495663
$ws.Range("O15").Value = 0.9636438974901603
$ws.Range("P15").Value = 0.9636438974901604
$ws.Range("Q15").Value = 2.8486086791
$ws.Range("R15").Value = 25.6374781119
$ws.Range("S15").Value = 0.04557345054142187
$ws.Range("T15").Value = 0.04557345054142188
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.41099
$ws.Range("H16").Value = 1.23297
$ws.Range("I16").Value = 0.0472928336495663
$ws.Range("J16").Value = 0.0472928336495663
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03527466666666667
$ws.Range("N16").Value = 0.105824
$ws.Range("O16").Value = 0.004904310471994002
$ws.Range("P16").Value = 0.004904310471994003
$ws.Range("Q16").Value = 0.01449753525333333
$ws.Range("R16").Value = 0.13047781728
$ws.Range("S16").Value = 0.0002319387393178383
$ws.Range("T16").Value = 0.0002319387393178384
$ws.Range("G17").Value = 3.718746666666667
$ws.Range("H17").Value = 11.15624
$ws.Range("I17").Value = 0.4279181184251341
$ws.Range("J17").Value = 0.4279181184251342
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.2262196666666666
$ws.Range("N17").Value = 0.6786589999999999
$ws.Range("O17").Value = 0.03145179203784564
$ws.Range("P17").Value = 0.03145179203784564
$ws.Range("Q17").Value = 0.841253631351111
$ws.Range("R17").Value = 7.571282682159999
$ws.Range("S17").Value = 0.01345879166993352
$ws.Range("T17").Value = 0.01345879166993352
$ws.Range("G18").Value = 3.718746666666667
$ws.Range("H18").Value = 11.15624
$ws.Range("I18").Value = 0.4279181184251341
$ws.Range("J18").Value = 0.4279181184251342
$ws.Range("O18").Value = 0.9636438974901603
$ws.Range("P18").Value = 0.9636438974901604
$ws.Range("Q18").Value = 25.77496783386667
$ws.Range("R18").Value = 231.9747105048
$ws.Range("S18").Value = 0.4123606834458522
$ws.Range("T18").Value = 0.4123606834458523
$ws.Range("G19").Value = 3.718746666666667
$ws.Range("H19").Value = 11.15624
$ws.Range("I19").Value = 0.4279181184251341
$ws.Range("J19").Value = 0.4279181184251342
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0.3333333333333333
$ws.Range("M19").Value = 0.03527466666666667
$ws.Range("N19").Value = 0.105824
$ws.Range("O19").Value = 0.004904310471994002
$ws.Range("P19").Value = 0.004904310471994003
$ws.Range("Q19").Value = 0.1311775490844445
$ws.Range("R19").Value = 1.18059794176
$ws.Range("S19").Value = 0.002098643309348355
$ws.Range("T19").Value = 0.002098643309348355
